$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Create the new "2022-Q4" sheet by duplicating "2022-Q3" (Item 2)
#    so it inherits identical header/column styling, then re-point
#    it right after "总计" (Item 1).
# ---------------------------------------------------------------
$templateSheet = $wb.Worksheets.Item(2)
$templateSheet.Copy($null, $wb.Worksheets.Item(1))
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Template sheet has 12 data rows (rows 2-13); we need 14 (rows 2-15).
# Extend column A styling (bold/border/center) down to the new rows.
$q4.Range("A13").Copy($q4.Range("A14:A15"))

# Helper: write a value into a cell while forcing TEXT storage (the
# source data keeps numeric-looking figures, like "12.39", as strings).
function Set-TextCell($ws, $row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$q4Data = @(
    @("012367", "上投摩根安荣回报混合C", "12.39", "25.70", "1.27", "0.1574", 7),
    @("001305", "九泰天富改革新动力混合A", "1.83", "93.65", "8.59", "0.1572", 6),
    @("004738", "上投摩根安隆回报混合A", "10.31", "23.32", "1.36", "0.1402", 5),
    @("012366", "上投摩根安荣回报混合A", "9.72", "25.70", "1.27", "0.1234", 7),
    @("001844", "九泰久益灵活配置混合C", "0.88", "94.08", "9.65", "0.0849", 6),
    @("004823", "上投摩根安裕回报混合A", "4.26", "36.12", "1.80", "0.0767", 7),
    @("004824", "上投摩根安裕回报混合C", "3.64", "36.12", "1.80", "0.0655", 7),
    @("004739", "上投摩根安隆回报混合C", "4.60", "23.32", "1.36", "0.0626", 5),
    @("001782", "九泰久益灵活配置混合A", "0.62", "94.08", "9.65", "0.0598", 6),
    @("004194", "招商中证1000指数增强A", "2.57", "94.27", "1.08", "0.0278", 6),
    @("005552", "富兰克林国海新趋势灵活配置混合A", "2.39", "20.31", "1.03", "0.0246", 5),
    @("004195", "招商中证1000指数增强C", "2.14", "94.27", "1.08", "0.0231", 6),
    @("005553", "富兰克林国海新趋势灵活配置混合C", "0.26", "20.31", "1.03", "0.0027", 5),
    @("009912", "九泰天富改革新动力混合C", "0.03", "93.65", "8.59", "0.0026", 6)
)

$r = 2
foreach ($row in $q4Data) {
    $q4.Cells.Item($r, 1).Value = ($r - 2)
    Set-TextCell $q4 $r 2 $row[0]
    Set-TextCell $q4 $r 3 $row[1]
    Set-TextCell $q4 $r 4 $row[2]
    Set-TextCell $q4 $r 5 $row[3]
    Set-TextCell $q4 $r 6 $row[4]
    Set-TextCell $q4 $r 7 $row[5]
    $q4.Cells.Item($r, 8).Value = $row[6]
    $r++
}

# ---------------------------------------------------------------
# 2) Update the "总计" summary sheet: insert the new 2022-Q4 row at
#    the top of the data and shift the existing rows down by one.
# ---------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

# Extend column A styling down to the new last row (row 8).
$total.Range("A7").Copy($total.Range("A8"))

$totalData = @(
    @("2022-Q4", 14, 1.01),
    @("2022-Q3", 12, 0.78),
    @("2022-Q2", 20, 1.52),
    @("2022-Q1", 12, 1.37),
    @("2021-Q4", 4, 0.27),
    @("2021-Q1", 12, 2.35),
    @("2020-Q4", 6, 2.25)
)

$r = 2
foreach ($row in $totalData) {
    $total.Cells.Item($r, 1).Value = ($r - 2)
    $total.Cells.Item($r, 2).Value = $row[0]
    $total.Cells.Item($r, 3).Value = $row[1]
    $total.Cells.Item($r, 4).Value = $row[2]
    $r++
}

Write-Host "2022-Q4 sheet added and 总计 summary updated"
